# Fruta / hortaliza, semanal
#
# A new weekly price observation for "Pepino ensalada" at "Terminal La
# Palmera de La Serena" needs to be inserted into the dataset. It belongs
# right after the existing row for date 44299 (row 511 in the 1-indexed
# sheet), pushing every subsequent record down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 511; Excel shifts rows 511..545 down to 512..546
# and copies the formatting (incl. the date style on column D) from the
# row above, exactly like a normal Excel "Insert Row" would.
$ws.Rows.Item(511).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A511").Value = 8
$ws.Range("B511").Value = "Terminal La Palmera de La Serena"
$ws.Range("C511").Value = "Coquimbo"
$ws.Range("D511").Value = 44714
$ws.Range("E511").Value = 4
$ws.Range("F511").Value = 100112043
$ws.Range("G511").Value = "Pepino ensalada"
$ws.Range("H511").Value = "Sin especificar"
$ws.Range("I511").Value = "Primera"
$ws.Range("J511").Value = 700
$ws.Range("K511").Value = 19000
$ws.Range("L511").Value = 20000
$ws.Range("M511").Value = 19500
$ws.Range("N511").Value = "`$/caja 60 unidades"
$ws.Range("O511").Value = "Región de Arica y Parinacota"
$ws.Range("P511").Value = 325
$ws.Range("Q511").Value = 60
$ws.Range("R511").Value = "Hortaliza"
